# edit.ps1 - apply the UNICEF Kenya sitrep "different calculations" update
#
# Summary of changes (see commit message "different calculations, need to
# make changes" / the accompanying OOXML diff):
#   1. Update the headline "children accessed basic education" figures and
#      county list.
#   2. Update the "teachers trained" figures.
#   3. Insert an editorial note (BlockText) plus a re-worded "EDC members
#      trained" paragraph (FirstParagraph) right before the old wording of
#      that same sentence, then delete the old (now-superseded) paragraph.
#   4. Add "Samburu" to the mentorship/life-skills county list.
#   5. Refresh every results figure in the "Sitrep table".
#   6. Append a new "Downloads" section (Heading2 + two hyperlinks) after
#      the table.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Headline paragraph: total children in basic education
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "A total of 4,290 children (2,100 girls and 2,190 boys) have accessed basic education with UNICEF support during May2023 in the counties of Isiolo, Wajir, Tana River, Kajiado, Kwale, Kilifi, Narok.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A total of 5,953 children (2,863 girls and 3,028 boys) have accessed basic education with UNICEF support during May2023 in the counties of Isiolo, Wajir, Tana River, Kajiado, Kwale, Kilifi, Narok, Baringo, Samburu.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Teachers trained figures
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "0 teachers (0 men and 0 women) were trained",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "168 teachers (77 men and 91 women) were trained",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3) Insert the editorial note + the re-worded "EDC members" paragraph
#    right before the existing (old-wording) paragraph, then delete the
#    old paragraph.
# ---------------------------------------------------------------------
$oldEdcPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*UNICEF also supported the training of 94 (67*") {
        $oldEdcPara = $p
        break
    }
}

if ($oldEdcPara -ne $null) {
    $insertPos = $oldEdcPara.Range.Start
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.InsertBefore(
        "This paragraph below has no clear match to any of the 5Ws activities. Not even when filtering for EDC-related activities do the figures tally.`r" +
        "UNICEF also supported the training of 94 members (67 men and 27 women) of the Enrolment Drive Committees (EDCs) and school Board of Management (BoM) in Baringo and Turkana counties on life skills, child safeguarding and their roles/responsibilities, child protection and referral mechanisms.`r"
    )

    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "This paragraph below has no clear match*") {
            $p.Style = "BlockText"
        }
        if ($p.Range.Text -like "UNICEF also supported the training of 94 members*") {
            $p.Style = "FirstParagraph"
        }
    }

    # Re-locate and remove the old (superseded) paragraph.
    $oldEdcPara2 = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*UNICEF also supported the training of 94 (67*") {
            $oldEdcPara2 = $p
            break
        }
    }
    if ($oldEdcPara2 -ne $null) {
        $oldEdcPara2.Range.Delete()
    }
}

# ---------------------------------------------------------------------
# 4) Add Samburu to the mentorship / life-skills county list
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Wajir, Tana River, Kajiado counties, strengthening",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Wajir, Tana River, Kajiado, Samburu counties, strengthening",
    2) | Out-Null

# ---------------------------------------------------------------------
# 5) Sitrep table: refresh every results figure
# ---------------------------------------------------------------------
$tableUpdates = @(
    @("38,262", "39,863"),
    @("4.41", "4.60"),
    @("148,686", "153,797"),
    @("113,400", "116,910"),
    @("12.38", "12.80"),
    @("18,260", "19,023"),
    @("4.21", "4.39"),
    @("67,020", "69,788"),
    @("51,697", "53,702"),
    @("11.16", "11.62"),
    @("20,002", "20,840"),
    @("4.62", "4.81"),
    @("81,666", "84,009"),
    @("61,703", "63,208"),
    @("13.60", "13.99"),
    @("134,664", "137,853"),
    @("22.42", "22.95"),
    @("63,137", "64,726"),
    @("21.03", "21.55"),
    @("72,327", "73,927"),
    @("24.09", "24.62")
)

foreach ($pair in $tableUpdates) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null
}

# ---------------------------------------------------------------------
# 6) Append the new "Downloads" section after the table
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$headingPara = $d.Paragraphs.Add($lastPara.Range)
$headingPara.Range.Text = "Downloads"
$headingPara.Style = "Heading2"

$lastPara2 = $d.Paragraphs.Last
$linkPara1 = $d.Paragraphs.Add($lastPara2.Range)
$linkPara1.Style = "FirstParagraph"
$linkPara1.Range.InsertAfter("5ws dataset")
$link1 = $d.Hyperlinks.Add($linkPara1.Range, "https://data.humdata.org/dataset/kenya-5ws-education", "", "", "5ws dataset")
$dup1Start = $link1.Range.End
$dup1End = $linkPara1.Range.End
if ($dup1End -gt $dup1Start) {
    $d.Range($dup1Start, $dup1End).Delete()
}

$lastPara3 = $d.Paragraphs.Last
$linkPara2 = $d.Paragraphs.Add($lastPara3.Range)
$linkPara2.Style = "BodyText"
$linkPara2.Range.InsertAfter("sitrep table")
$link2 = $d.Hyperlinks.Add($linkPara2.Range, "https://data.humdata.org/dataset/kenya-education-sitrep-table", "", "", "sitrep table")
$dup2Start = $link2.Range.End
$dup2End = $linkPara2.Range.End
if ($dup2End -gt $dup2Start) {
    $d.Range($dup2Start, $dup2End).Delete()
}

Write-Output "edit complete"
